$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.652167666666667
$ws.Range("H2").Value = 16.956503
$ws.Range("I2").Value = 0.1860329065948871
$ws.Range("J2").Value = 0.1860329065948871
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.652835333333333
$ws.Range("N2").Value = 10.958506
$ws.Range("O2").Value = 0.1451640905049821
$ws.Range("P2").Value = 0.1451640905049821
$ws.Range("Q2").Value = 20.64643776272423
$ws.Range("R2").Value = 185.817939864518
$ws.Range("S2").Value = 0.02700529768984507
$ws.Range("T2").Value = 0.02700529768984507

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.652167666666667
$ws.Range("H3").Value = 16.956503
$ws.Range("I3").Value = 0.1860329065948871
$ws.Range("J3").Value = 0.1860329065948871
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.483777
$ws.Range("N3").Value = 7.451331
$ws.Range("O3").Value = 0.09870557972652284
$ws.Range("P3").Value = 0.09870557972652286
$ws.Range("Q3").Value = 14.03872405061033
$ws.Range("R3").Value = 126.348516455493
$ws.Range("S3").Value = 0.0183624858936584
$ws.Range("T3").Value = 0.01836248589365841

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.652167666666667
$ws.Range("H4").Value = 16.956503
$ws.Range("I4").Value = 0.1860329065948871
$ws.Range("J4").Value = 0.1860329065948871
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 19.026879
$ws.Range("N4").Value = 57.080637
$ws.Range("O4").Value = 0.7561303297684949
$ws.Range("P4").Value = 0.756130329768495
$ws.Range("Q4").Value = 107.543110281379
$ws.Range("R4").Value = 967.887992532411
$ws.Range("S4").Value = 0.1406651230113836
$ws.Range("T4").Value = 0.1406651230113836

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.407289
$ws.Range("H5").Value = 49.221867
$ws.Range("I5").Value = 0.5400221369958743
$ws.Range("J5").Value = 0.5400221369958743
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.652835333333333
$ws.Range("N5").Value = 10.958506
$ws.Range("O5").Value = 0.1451640905049821
$ws.Range("P5").Value = 0.1451640905049821
$ws.Range("Q5").Value = 59.93312498341135
$ws.Range("R5").Value = 539.3981248507021
$ws.Range("S5").Value = 0.07839182236956296
$ws.Range("T5").Value = 0.07839182236956296

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 16.407289
$ws.Range("H6").Value = 49.221867
$ws.Range("I6").Value = 0.5400221369958743
$ws.Range("J6").Value = 0.5400221369958743
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.483777
$ws.Range("N6").Value = 7.451331
$ws.Range("O6").Value = 0.09870557972652284
$ws.Range("P6").Value = 0.09870557972652286
$ws.Range("Q6").Value = 40.752047050553
$ws.Range("R6").Value = 366.768423454977
$ws.Range("S6").Value = 0.05330319809733351
$ws.Range("T6").Value = 0.05330319809733352

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 16.407289
$ws.Range("H7").Value = 49.221867
$ws.Range("I7").Value = 0.5400221369958743
$ws.Range("J7").Value = 0.5400221369958743
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 19.026879
$ws.Range("N7").Value = 57.080637
$ws.Range("O7").Value = 0.7561303297684949
$ws.Range("P7").Value = 0.756130329768495
$ws.Range("Q7").Value = 312.179502521031
$ws.Range("R7").Value = 2809.615522689279
$ws.Range("S7").Value = 0.4083271165289778
$ws.Range("T7").Value = 0.4083271165289779

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.323166333333333
$ws.Range("H8").Value = 24.969499
$ws.Range("I8").Value = 0.2739449564092387
$ws.Range("J8").Value = 0.2739449564092387
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.652835333333333
$ws.Range("N8").Value = 10.958506
$ws.Range("O8").Value = 0.1451640905049821
$ws.Range("P8").Value = 0.1451640905049821
$ws.Range("Q8").Value = 30.40315606761044
$ws.Range("R8").Value = 273.628404608494
$ws.Range("S8").Value = 0.0397669704455741
$ws.Range("T8").Value = 0.03976697044557412

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.323166333333333
$ws.Range("H9").Value = 24.969499
$ws.Range("I9").Value = 0.2739449564092387
$ws.Range("J9").Value = 0.2739449564092387
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.483777
$ws.Range("N9").Value = 7.451331
$ws.Range("O9").Value = 0.09870557972652284
$ws.Range("P9").Value = 0.09870557972652286
$ws.Range("Q9").Value = 20.67288910590766
$ws.Range("R9").Value = 186.056001953169
$ws.Range("S9").Value = 0.02703989573553093
$ws.Range("T9").Value = 0.02703989573553094

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.323166333333333
$ws.Range("H10").Value = 24.969499
$ws.Range("I10").Value = 0.2739449564092387
$ws.Range("J10").Value = 0.2739449564092387
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 19.026879
$ws.Range("N10").Value = 57.080637
$ws.Range("O10").Value = 0.7561303297684949
$ws.Range("P10").Value = 0.756130329768495
$ws.Range("Q10").Value = 158.363878721207
$ws.Range("R10").Value = 1425.274908490863
$ws.Range("S10").Value = 0.2071380902281336
$ws.Range("T10").Value = 0.2071380902281337
